$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ValidateFormulas")

# Insert a new row before row 31 so that the old row 31 content shifts down to row 32,
# and everything below shifts down by one as well.
$ws.Rows.Item(31).Insert()

# Populate the newly inserted row 31 with the "If" example data.
$ws.Range("A31").Value = "If"
$ws.Range("B31").Formula = "=IF(B2>3,B3,B5)"
$ws.Range("C31").Formula = "=IF((B2*B3)*C1<0,(B2*B3)*C1,ABS((B2*B3)*C1))"
$ws.Range("D31").Formula = "=IF((B2*B3)*C1<0,ABS((B2*B3)*C1),(B2*B3)*C1)"

# Apply the same bold formatting as other "label" cells in column A.
$ws.Range("A31").Font.Bold = $true

# Update the selected cell to D31 as shown in the diff.
$ws.Range("D31").Select()

$wb.Save()
